# PROS-6288 - MOLSONCOORSHR - development
# Rename KPI "Linear SOS vs Target" -> "Facings SOS vs Target" wherever it
# appears on the KPIs sheet, and update the sheet's selection/cursor state
# to match the saved workbook view (final selection resting on E9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "Linear SOS vs Target"
$newText = "Facings SOS vs Target"

$usedRange = $ws.UsedRange
$rows = $usedRange.Rows.Count
$cols = $usedRange.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Value2 -eq $oldText) {
            $cell.Value = $newText
        }
    }
}

# Move the on-screen selection: header/frozen pane back to the first column,
# and the scrollable data pane down to the last data row (E9), matching the
# workbook's saved view state.
$ws.Range("A1").Select()
$ws.Range("E9").Select()
